$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rng, $val) {
    $oldFormat = $rng.NumberFormat
    $oldStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = $oldFormat
    $rng.Style = $oldStyle
}

$ws.Range("D2").Value = '30.386.06'
$ws.Range("E2").Value = '  +0.15%  '
$ws.Range("D3").Value = '1.876.32'
$ws.Range("E3").Value = '  -0.62%  '
$ws.Range("E4").Value = '  -0.08%  '
Set-TextValue $ws.Range("D5") '238.92'
$ws.Range("E5").Value = '  +0.43%  '
Set-TextValue $ws.Range("D6") '1.001'
Set-TextValue $ws.Range("D7") '0.4799'
$ws.Range("E7").Value = '  -0.35%  '
Set-TextValue $ws.Range("D8") '0.2820'
$ws.Range("E8").Value = '  -2.62%  '
Set-TextValue $ws.Range("D9") '0.06523'
$ws.Range("E9").Value = '  -1.26%  '
$ws.Range("D10").Value = '1.874.92'
$ws.Range("E10").Value = '  -1.26%  '
Set-TextValue $ws.Range("D11") '0.07484'
$ws.Range("E11").Value = '  +1.39%  '
Set-TextValue $ws.Range("D12") '16.57'
$ws.Range("E12").Value = '  -2.07%  '
Set-TextValue $ws.Range("D13") '5.070'
$ws.Range("E13").Value = '  -2.02%  '
Set-TextValue $ws.Range("D14") '88.41'
$ws.Range("E14").Value = '  +0.71%  '
Set-TextValue $ws.Range("D15") '0.6613'
$ws.Range("E15").Value = '  -0.24%  '
$ws.Range("D16").Value = '30.351.23'
$ws.Range("E16").Value = '  +0.08%  '
Set-TextValue $ws.Range("D17") '13.30'
$ws.Range("E17").Value = '  -1.27%  '
Set-TextValue $ws.Range("D18") '0.9997'
$ws.Range("E18").Value = '  -0.06%  '
Set-TextValue $ws.Range("D19") '0.000007600'
$ws.Range("E19").Value = '  -2.03%  '
$ws.Range("D20").Value = '2.115.49'
$ws.Range("E20").Value = '  -1.53%  '
Set-TextValue $ws.Range("D21") '5.304'
$ws.Range("E21").Value = '  -2.84%  '
Set-TextValue $ws.Range("D22") '0.9997'
$ws.Range("E22").Value = '  -0.07%  '
Set-TextValue $ws.Range("D23") '220.85'
$ws.Range("E23").Value = '  +15.21%  '
Set-TextValue $ws.Range("D24") '6.194'
$ws.Range("E24").Value = '  -0.21%  '
Set-TextValue $ws.Range("D25") '9.340'
$ws.Range("E25").Value = '  -1.18%  '
Set-TextValue $ws.Range("D26") '167.57'
$ws.Range("E27").Value = '  +0.81%  '
Set-TextValue $ws.Range("D28") '1.961'
$ws.Range("E28").Value = '  +1.12%  '
Set-TextValue $ws.Range("D29") '1.456'
$ws.Range("E29").Value = '  +0.30%  '
$ws.Range("E30").Value = '  +2.40%  '
Set-TextValue $ws.Range("D31") '4.306'
$ws.Range("E31").Value = '  +1.04%  '
Set-TextValue $ws.Range("D32") '4.025'
$ws.Range("E32").Value = '  -0.75%  '
Set-TextValue $ws.Range("D33") '0.05011'
$ws.Range("E33").Value = '  -1.31%  '
Set-TextValue $ws.Range("D34") '1.210'
$ws.Range("E34").Value = '  +5.83%  '
Set-TextValue $ws.Range("D35") '0.7433'
$ws.Range("E35").Value = '  +0.46%  '
Set-TextValue $ws.Range("D36") '2.705'
$ws.Range("E36").Value = '  -0.36%  '
Set-TextValue $ws.Range("D37") '0.01826'
$ws.Range("E37").Value = '  -0.13%  '
Set-TextValue $ws.Range("D38") '2.615'
$ws.Range("E38").Value = '  -1.25%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range("D39") '2.063'
$ws.Range("E39").Value = '  -1.26%  '
$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range("D40") '0.9050'
$ws.Range("E40").Value = '  -1.69%  '
Set-TextValue $ws.Range("D41") '106.39'
$ws.Range("E41").Value = '  +0.15%  '
Set-TextValue $ws.Range("D42") '5.856'
$ws.Range("E42").Value = '  -0.92%  '
Set-TextValue $ws.Range("D43") '0.4275'
$ws.Range("E43").Value = '  -1.36%  '
$ws.Range("E44").Value = '  +0.34%  '
Set-TextValue $ws.Range("D45") '7.434'
$ws.Range("E45").Value = '  -2.87%  '
Set-TextValue $ws.Range("D46") '64.37'
$ws.Range("E46").Value = '  -1.04%  '
Set-TextValue $ws.Range("D47") '0.1275'
$ws.Range("E47").Value = '  -7.22%  '
Set-TextValue $ws.Range("D48") '1.474'
$ws.Range("E48").Value = '  -7.02%  '
Set-TextValue $ws.Range("D49") '8.866'
$ws.Range("E49").Value = '  -1.61%  '
Set-TextValue $ws.Range("D50") '33.70'
$ws.Range("E50").Value = '  -1.75%  '
Set-TextValue $ws.Range("D51") '0.3886'
$ws.Range("E51").Value = '  +0.28%  '
